$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 'Dinesh Chandimal'
$ws.Range("C2").Value = 4
$ws.Range("D2").Value = 'LBW'
$ws.Range("E2").Value = ' Josh Hazlewood'
$ws.Range("J2").Value = 'David Warner'
$ws.Range("M2").Value = 'LBW'
$ws.Range("N2").Value = ' Nuwan Pradeep'
# Row 3
$ws.Range("A3").Value = 'Pathum Nissanka'
$ws.Range("B3").Value = 0
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = 'LBW'
$ws.Range("E3").Value = ' Adam Zampa'
$ws.Range("J3").Value = 'Aaron Finch(C)'
$ws.Range("K3").Value = 5
$ws.Range("M3").Value = 'LBW'
$ws.Range("N3").Value = ' Nuwan Pradeep'
# Row 4
$ws.Range("A4").Value = 'Charith Asalanka'
$ws.Range("B4").Value = 90
$ws.Range("C4").Value = 35
$ws.Range("D4").Value = 'Bowled'
$ws.Range("E4").Value = ' Mitchell Starc'
$ws.Range("J4").Value = 'Mitchell Marsh'
$ws.Range("K4").Value = 19
$ws.Range("M4").Value = 'LBW'
$ws.Range("N4").Value = ' Wanindu Hasaranga'
# Row 5
$ws.Range("A5").Value = 'Dhananjaya de Silva'
$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 1
$ws.Range("E5").Value = ' Adam Zampa'
$ws.Range("J5").Value = 'Steve Smith'
$ws.Range("K5").Value = 21
$ws.Range("L5").Value = 7
$ws.Range("M5").Value = 'Bowled'
$ws.Range("N5").Value = ' Dushmantha Chameera'
# Row 6
$ws.Range("A6").Value = 'Bhanuka Rajapakse'
$ws.Range("B6").Value = 21
$ws.Range("C6").Value = 10
$ws.Range("D6").Value = 'LBW'
$ws.Range("E6").Value = ' Josh Hazlewood'
$ws.Range("J6").Value = 'Glenn Maxwell'
$ws.Range("K6").Value = 57
$ws.Range("L6").Value = 22
$ws.Range("M6").Value = 'Bowled'
$ws.Range("N6").Value = ' Dushmantha Chameera'
# Row 7
$ws.Range("A7").Value = 'Dasun Shanka(C)'
$ws.Range("B7").Value = 6
$ws.Range("C7").Value = 2
$ws.Range("E7").Value = ' Adam Zampa'
$ws.Range("J7").Value = 'Matthew Wade'
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 2
$ws.Range("M7").Value = 'Bowled'
$ws.Range("N7").Value = ' Nuwan Pradeep'
# Row 8
$ws.Range("A8").Value = 'Wanindu Hasaranga'
$ws.Range("B8").Value = 37
$ws.Range("C8").Value = 13
$ws.Range("E8").Value = ' Marcus Stionis'
$ws.Range("J8").Value = 'Marcus Stionis'
$ws.Range("K8").Value = 17
$ws.Range("L8").Value = 8
$ws.Range("M8").Value = 'Bowled'
$ws.Range("N8").Value = ' Wanindu Hasaranga'
# Row 9
$ws.Range("A9").Value = 'Chamika Karunarathne'
$ws.Range("C9").Value = 3
$ws.Range("D9").Value = 'Caught'
$ws.Range("E9").Value = ' Marcus Stionis'
$ws.Range("J9").Value = 'Pat Cummins'
$ws.Range("K9").Value = 14
$ws.Range("L9").Value = 11
$ws.Range("M9").Value = 'Caught'
$ws.Range("N9").Value = ' Nuwan Pradeep'
# Row 10
$ws.Range("A10").Value = 'Dushmantha Chameera'
$ws.Range("B10").Value = 4
$ws.Range("C10").Value = 2
$ws.Range("D10").Value = 'Caught'
$ws.Range("E10").Value = ' Josh Hazlewood'
$ws.Range("J10").Value = 'Mitchell Starc'
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 'Caught'
$ws.Range("N10").Value = ' Wanindu Hasaranga'
# Row 11
$ws.Range("A11").Value = 'Maheesh Theekshana'
$ws.Range("B11").Value = 23
$ws.Range("C11").Value = 6
$ws.Range("D11").Value = 'NOT OUT'
$ws.Range("E11").Value = ' '
$ws.Range("J11").Value = 'Adam Zampa'
$ws.Range("K11").Value = 25
$ws.Range("L11").Value = 11
$ws.Range("M11").Value = 'NOT OUT'
$ws.Range("N11").Value = ' '
# Row 12
$ws.Range("A12").Value = 'Nuwan Pradeep'
$ws.Range("B12").Value = 14
$ws.Range("C12").Value = 8
$ws.Range("E12").Value = ' Marcus Stionis'
$ws.Range("J12").Value = 'Josh Hazlewood'
$ws.Range("K12").Value = 21
$ws.Range("L12").Value = 11
$ws.Range("M12").Value = 'LBW'
$ws.Range("N12").Value = ' Chamika Karunarathne'
# Row 16 (innings totals)
$ws.Range("A16").Value = 211
$ws.Range("C16").Value = "'14.1"
$ws.Range("C16").Style = "Normal"
$ws.Range("D16").Value = 85
$ws.Range("J16").Value = 179
$ws.Range("L16").Value = "'14.0"
$ws.Range("L16").Style = "Normal"
$ws.Range("M16").Value = 84
# Row 21 (bowling figures)
$ws.Range("A21").Value = 'Josh Hazlewood'
$ws.Range("B21").Value = "'3.0"
$ws.Range("B21").Style = "Normal"
$ws.Range("C21").Value = 55
$ws.Range("D21").Value = 3
$ws.Range("E21").Value = 18.33
$ws.Range("J21").Value = 'Wanindu Hasaranga'
$ws.Range("K21").Value = "'2.0"
$ws.Range("K21").Style = "Normal"
$ws.Range("L21").Value = 21
$ws.Range("N21").Value = 10.5
# Row 22
$ws.Range("A22").Value = 'Adam Zampa'
$ws.Range("B22").Value = "'3.0"
$ws.Range("B22").Style = "Normal"
$ws.Range("C22").Value = 41
$ws.Range("D22").Value = 3
$ws.Range("E22").Value = 13.67
$ws.Range("J22").Value = 'Nuwan Pradeep'
$ws.Range("K22").Value = "'3.0"
$ws.Range("K22").Style = "Normal"
$ws.Range("L22").Value = 29
$ws.Range("M22").Value = 4
$ws.Range("N22").Value = 9.67
# Row 23
$ws.Range("A23").Value = 'Mitchell Starc'
$ws.Range("B23").Value = "'3.0"
$ws.Range("B23").Style = "Normal"
$ws.Range("C23").Value = 47
$ws.Range("D23").Value = 1
$ws.Range("E23").Value = 15.67
$ws.Range("J23").Value = 'Maheesh Theekshana'
$ws.Range("K23").Value = "'3.0"
$ws.Range("K23").Style = "Normal"
$ws.Range("L23").Value = 51
$ws.Range("N23").Value = 17
# Row 24
$ws.Range("A24").Value = 'Pat Cummins'
$ws.Range("C24").Value = 36
$ws.Range("D24").Value = 0
$ws.Range("E24").Value = 12
$ws.Range("J24").Value = 'Dushmantha Chameera'
$ws.Range("K24").Value = "'3.0"
$ws.Range("K24").Style = "Normal"
$ws.Range("L24").Value = 44
$ws.Range("M24").Value = 2
$ws.Range("N24").Value = 14.67
# Row 25
$ws.Range("A25").Value = 'Marcus Stionis'
$ws.Range("B25").Value = "'2.1"
$ws.Range("B25").Style = "Normal"
$ws.Range("C25").Value = 32
$ws.Range("D25").Value = 3
$ws.Range("E25").Value = 15.24
$ws.Range("J25").Value = 'Chamika Karunarathne'
$ws.Range("L25").Value = 34
$ws.Range("M25").Value = 1
$ws.Range("N25").Value = 11.33
